$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("D2").Value = "'69.362.35"
$ws.Range("E2").Value = '  +1.73%  '

# Row 3
$ws.Range("B3").Value = 'Ethereum'
$ws.Range("D3").Value = "'3.901.54"
$ws.Range("E3").Value = '  -0.22%  '

# Row 4
$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("B5").Value = 'BNB'
$ws.Range("D5").Value = "'527.63"
$ws.Range("E5").Value = '  +9.13%  '

# Row 6
$ws.Range("B6").Value = 'Solana'
$ws.Range("D6").Value = "'143.51"
$ws.Range("E6").Value = '  -1.96%  '

# Row 7
$ws.Range("B7").Value = 'XRP'
$ws.Range("D7").Value = "'0.610"
$ws.Range("E7").Value = '  -1.67%  '

# Row 8
$ws.Range("B8").Value = 'USDC'
$ws.Range("D8").Value = "'0.998"
$ws.Range("E8").Value = '  +0.05%  '

# Row 9
$ws.Range("B9").Value = 'Cardano'
$ws.Range("D9").Value = "'0.721"
$ws.Range("E9").Value = '  -1.59%  '

# Row 10
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("D10").Value = "'0.171"
$ws.Range("E10").Value = '  +1.96%  '

# Row 11
$ws.Range("B11").Value = 'ShibaInu'
$ws.Range("D11").Value = "'0.0000333"
$ws.Range("E11").Value = '  -4.02%  '

# Row 12
$ws.Range("B12").Value = 'Avalanche'
$ws.Range("D12").Value = "'42.02"
$ws.Range("E12").Value = '  -2.48%  '

# Row 13
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("D13").Value = "'4.516.68"
$ws.Range("E13").Value = '  -0.45%  '

# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("D14").Value = "'10.21"
$ws.Range("E14").Value = '  -4.38%  '

# Row 15
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("D15").Value = "'3.897.89"
$ws.Range("E15").Value = '  -0.91%  '

# Row 16
$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").Value = "'1.23"
$ws.Range("E16").Value = '  +8.38%  '

# Row 17
$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").Value = "'0.135"
$ws.Range("E17").Value = '  -0.60%  '

# Row 18
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("D18").Value = "'13.80"
$ws.Range("E18").Value = '  -2.79%  '

# Row 19
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("D19").Value = "'19.70"
$ws.Range("E19").Value = '  -2.33%  '

# Row 20
$ws.Range("B20").Value = 'WrappedBTC'
$ws.Range("D20").Value = "'69.251.48"
$ws.Range("E20").Value = '  +1.52%  '

# Row 21
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("D21").Value = "'426.22"
$ws.Range("E21").Value = '  -0.70%  '

# Row 22
$ws.Range("B22").Value = 'ImmutableX'
$ws.Range("D22").Value = "'3.34"
$ws.Range("E22").Value = '  -4.49%  '

# Row 23
$ws.Range("B23").Value = 'InternetComputer(DFINITY)'
$ws.Range("D23").Value = "'14.16"
$ws.Range("E23").Value = '  -6.44%  '

# Row 24
$ws.Range("B24").Value = 'PancakeSwap'
$ws.Range("D24").Value = "'4.09"
$ws.Range("E24").Value = '  +10.67%  '

# Row 25
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("D25").Value = "'87.66"
$ws.Range("E25").Value = '  -0.72%  '

# Row 26
$ws.Range("B26").Value = 'RenderToken'
$ws.Range("D26").Value = "'11.61"
$ws.Range("E26").Value = '  -0.57%  '

# Row 27
$ws.Range("B27").Value = 'Filecoin'
$ws.Range("D27").Value = "'10.55"
$ws.Range("E27").Value = '  -5.39%  '

# Row 28
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("D28").Value = "'35.95"
$ws.Range("E28").Value = '  -4.66%  '

# Row 29
$ws.Range("B29").Value = 'Bittensor'
$ws.Range("D29").Value = "'693.50"
$ws.Range("E29").Value = '  -3.23%  '

# Row 30
$ws.Range("B30").Value = 'Cosmos'
$ws.Range("D30").Value = "'13.13"
$ws.Range("E30").Value = '  -5.03%  '

# Row 31
$ws.Range("B31").Value = 'Hedera'
$ws.Range("D31").Value = "'0.125"
$ws.Range("E31").Value = '  -4.11%  '

# Row 32
$ws.Range("B32").Value = 'Toncoin'
$ws.Range("D32").Value = "'2.82"
$ws.Range("E32").Value = '  -3.55%  '

# Row 33
$ws.Range("B33").Value = 'OKB'
$ws.Range("D33").Value = "'67.65"
$ws.Range("E33").Value = '  +11.94%  '

# Row 34
$ws.Range("B34").Value = 'TheGraph'
$ws.Range("D34").Value = "'0.443"
$ws.Range("E34").Value = '  +11.98%  '

# Row 35
$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("D35").Value = "'5.95"
$ws.Range("E35").Value = '  -4.17%  '

# Row 36
$ws.Range("B36").Value = 'InjectiveProtocol'
$ws.Range("D36").Value = "'40.23"
$ws.Range("E36").Value = '  -3.18%  '

# Row 37
$ws.Range("B37").Value = 'PEPE'
$ws.Range("D37").Value = "'0.0₃0840"
$ws.Range("E37").Value = '  -7.20%  '

# Row 38
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").Value = "'0.148"
$ws.Range("E38").Value = '  +4.17%  '

# Row 39
$ws.Range("B39").Value = 'Dai'
$ws.Range("C39").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D39").Value = "'0.999"
$ws.Range("E39").Value = '  +0.06%  '

# Row 40
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = '  -0.12%  '

# Row 41
$ws.Range("B41").Value = 'VeChain'
$ws.Range("D41").Value = "'0.0479"
$ws.Range("E41").Value = '  -1.92%  '

# Row 42
$ws.Range("B42").Value = 'Fetch.AI'
$ws.Range("D42").Value = "'2.79"
$ws.Range("E42").Value = '  -6.81%  '

# Row 43
$ws.Range("B43").Value = 'WEMIXToken'
$ws.Range("D43").Value = "'3.00"
$ws.Range("E43").Value = '  -0.12%  '

# Row 44
$ws.Range("B44").Value = 'ThetaToken'
$ws.Range("D44").Value = "'2.94"
$ws.Range("E44").Value = '  -5.87%  '

# Row 45
$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("D45").Value = "'3.34"
$ws.Range("E45").Value = '  +0.25%  '

# Row 46
$ws.Range("B46").Value = 'Stellar'
$ws.Range("D46").Value = "'0.140"
$ws.Range("E46").Value = '  -1.65%  '

# Row 47
$ws.Range("B47").Value = 'Stacks'
$ws.Range("D47").Value = "'3.06"
$ws.Range("E47").Value = '  +9.17%  '

# Row 48
$ws.Range("B48").Value = 'LidoDAOToken'
$ws.Range("D48").Value = "'3.28"
$ws.Range("E48").Value = '  -4.73%  '

# Row 49
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = "'0.0₆0337"
$ws.Range("E49").Value = '  +0.07%  '

# Row 50
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = "'142.47"
$ws.Range("E50").Value = '  -1.51%  '

# Row 51
$ws.Range("B51").Value = 'ARBITRUM'
$ws.Range("C51").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D51").Value = "'2.05"
$ws.Range("E51").Value = '  -3.84%  '
